# Fix Digital Album Sales to correctly match previous royalty runs.
# Updates unit-sales counts (column I) on the "RS Digital Sales Template v1_23"
# sheet for the rows whose figures were recalculated, then leaves the
# selection where the author last left it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I2").Value = 23
$ws.Range("I3").Value = 9
$ws.Range("I4").Value = 29
$ws.Range("I5").Value = 22
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 20
$ws.Range("I8").Value = 31
$ws.Range("I10").Value = 80
$ws.Range("I11").Value = 15
$ws.Range("I14").Value = 13
$ws.Range("I16").Value = 6
$ws.Range("I17").Value = 6
$ws.Range("I18").Value = 2
$ws.Range("I19").Value = 11
$ws.Range("I20").Value = 3
$ws.Range("I21").Value = 2
$ws.Range("I22").Value = 5
$ws.Range("I23").Value = 7
$ws.Range("I24").Value = 5
$ws.Range("I25").Value = 5
$ws.Range("I26").Value = 14
$ws.Range("I31").Value = 35
$ws.Range("I32").Value = 8
$ws.Range("I33").Value = 37
$ws.Range("I34").Value = 5
$ws.Range("I35").Value = 3
$ws.Range("I36").Value = 8
$ws.Range("I37").Value = 19
$ws.Range("I38").Value = 19
$ws.Range("I39").Value = 16
$ws.Range("I40").Value = 12
$ws.Range("I41").Value = 12
$ws.Range("I42").Value = 10
$ws.Range("I43").Value = 15
$ws.Range("I44").Value = 11
$ws.Range("I45").Value = 38
$ws.Range("I46").Value = 40
$ws.Range("I49").Value = 9
$ws.Range("I50").Value = 16
$ws.Range("I51").Value = 80
$ws.Range("I52").Value = 31
$ws.Range("I54").Value = 0
$ws.Range("I56").Value = 4
$ws.Range("I58").Value = 9
$ws.Range("I60").Value = 91
$ws.Range("I62").Value = 44
$ws.Range("I64").Value = 8
$ws.Range("I65").Value = 11
$ws.Range("I67").Value = 7
$ws.Range("I69").Value = 26
$ws.Range("I70").Value = 43
$ws.Range("I72").Value = 5
$ws.Range("I73").Value = 64
$ws.Range("M10").Select()
